# Apply updated crypto price/volume data per diff (cryptos list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.735.84"
$ws.Range("E2").Value = "  -1.70%  "
$ws.Range("D3").Value = "'2.890.05"
$ws.Range("E3").Value = "  -1.77%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'565.58"
$ws.Range("E5").Value = "  -4.44%  "
$ws.Range("D6").Value = "'142.93"
$ws.Range("E6").Value = "  -2.89%  "
$ws.Range("D8").Value = "'0.501"
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("D9").Value = "'2.888.50"
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("D10").Value = "'6.93"
$ws.Range("E10").Value = "  -1.76%  "
$ws.Range("D11").Value = "'0.146"
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("D12").Value = "'0.428"
$ws.Range("E12").Value = "  -2.07%  "
$ws.Range("E13").Value = "  -0.73%  "
$ws.Range("D14").Value = "'31.66"
$ws.Range("E14").Value = "  -2.36%  "
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").Value = "'3.370.60"
$ws.Range("E16").Value = "  -1.71%  "
$ws.Range("D17").Value = "'61.679.54"
$ws.Range("E17").Value = "  -1.70%  "
$ws.Range("D18").Value = "'2.900.99"
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("E19").Value = "  -2.44%  "
$ws.Range("D20").Value = "'429.60"
$ws.Range("E20").Value = "  -1.87%  "
$ws.Range("D21").Value = "'13.01"
$ws.Range("E21").Value = "  -2.84%  "
$ws.Range("E22").Value = "  -1.97%  "
$ws.Range("D23").Value = "'6.79"
$ws.Range("E23").Value = "  -2.83%  "
$ws.Range("D24").Value = "'78.93"
$ws.Range("E24").Value = "  -2.19%  "
$ws.Range("D25").Value = "'11.89"
$ws.Range("E25").Value = "  +0.83%  "
$ws.Range("D27").Value = "'9.92"
$ws.Range("E27").Value = "  -11.19%  "
$ws.Range("E28").Value = "  -5.25%  "
$ws.Range("E29").Value = "  +8.69%  "
$ws.Range("D30").Value = "'7.01"
$ws.Range("E30").Value = "  -3.56%  "
$ws.Range("E31").Value = "  -4.39%  "
$ws.Range("D32").Value = "'2.03"
$ws.Range("E32").Value = "  -8.56%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("E34").Value = "  -1.44%  "
$ws.Range("E35").Value = "  -3.30%  "
$ws.Range("E36").Value = "  -3.18%  "
$ws.Range("D37").Value = "'5.36"
$ws.Range("E37").Value = "  -4.19%  "
$ws.Range("D38").Value = "'48.81"
$ws.Range("E38").Value = "  -1.51%  "
$ws.Range("E39").Value = "  -4.64%  "
$ws.Range("D40").Value = "'2.80"
$ws.Range("E40").Value = "  -6.62%  "
$ws.Range("E41").Value = "  -3.30%  "
$ws.Range("E42").Value = "  -3.69%  "
$ws.Range("D43").Value = "'39.40"
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("E44").Value = "  -4.14%  "
$ws.Range("D45").Value = "'2.679.11"
$ws.Range("E45").Value = "  -0.75%  "
$ws.Range("D46").Value = "'132.14"
$ws.Range("E46").Value = "  -2.26%  "
$ws.Range("D47").Value = "'0.0335"
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("D48").Value = "'343.10"
$ws.Range("E48").Value = "  -3.82%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  -1.83%  "
$ws.Range("D51").Value = "'21.51"
$ws.Range("E51").Value = "  -4.86%  "
